# Add a new "hole size" entry in row 9 (L9:P9), mirroring the pattern
# established by the rows above it (which use shared formulas across
# M3:M8 / N3:N8 / O3:O8 / P3:P8).
#
# L9 is the new raw "intended hole size" input (2.9); M9:P9 are the same
# formulas used in the rows above, re-pointed at row 9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New input value
$ws.Range("L9").Value = 2.9

# Formulas, following the same pattern as M3/N3/O3/P3 (etc.), just for row 9
$ws.Range("M9").Formula = "=(L9+0.2926) / 0.9917"
$ws.Range("N9").Formula = "=L9+H9"
$ws.Range("O9").Formula = "=ROUND(MAX(M9:N9), 1)"
$ws.Range("P9").Formula = "=O9/2"

# Recalculate so cached formula results are written out
$wb.Application.Calculate()

# Reflect the cursor/selection position left behind in the saved file
$ws.Range("M13").Select()
